{"js": "// Replace the 25 division problems in the table with their new values.\n// Each old string is unique in the document, so we resolve every\n// search() hit against the *original* text first, then perform the\n// text replacements. Doing the lookups before any mutation avoids any\n// chance of a later replacement accidentally re-matching text that was\n// just written by an earlier replacement (a couple of the new values\n// equal other cells' old values).\nconst replacements = [\n  [\"34\u00f74=\", \"24\u00f79=\"],\n  [\"56\u00f73=\", \"18\u00f73=\"],\n  [\"43\u00f72=\", \"49\u00f73=\"],\n  [\"95\u00f75=\", \"26\u00f79=\"],\n  [\"52\u00f75=\", \"99\u00f73=\"],\n  [\"48\u00f75=\", \"32\u00f76=\"],\n  [\"12\u00f75=\", \"55\u00f72=\"],\n  [\"20\u00f73=\", \"62\u00f74=\"],\n  [\"14\u00f74=\", \"94\u00f79=\"],\n  [\"65\u00f74=\", \"54\u00f73=\"],\n  [\"90\u00f77=\", \"96\u00f76=\"],\n  [\"57\u00f77=\", \"81\u00f77=\"],\n  [\"80\u00f73=\", \"43\u00f73=\"],\n  [\"29\u00f75=\", \"57\u00f77=\"],\n  [\"45\u00f77=\", \"69\u00f72=\"],\n  [\"19\u00f72=\", \"13\u00f73=\"],\n  [\"45\u00f79=\", \"38\u00f74=\"],\n  [\"59\u00f77=\", \"17\u00f79=\"],\n  [\"26\u00f79=\", \"93\u00f77=\"],\n  [\"27\u00f77=\", \"49\u00f75=\"],\n  [\"10\u00f74=\", \"51\u00f76=\"],\n  [\"74\u00f72=\", \"99\u00f74=\"],\n  [\"15\u00f76=\", \"23\u00f72=\"],\n  [\"26\u00f76=\", \"77\u00f77=\"],\n  [\"57\u00f73=\", \"36\u00f74=\"],\n];\n\nconst searchResults = [];\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  searchResults.push({ newText, found, oldText });\n}\nawait context.sync();\n\nfor (const { newText, found, oldText } of searchResults) {\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text \"${oldText}\" to replace.`);\n  }\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 25 division problems in the worksheet's table with their\n# new values. The problems live in table 1, rows 1/5/9/13/17 (the other\n# rows are blank \"show your work\" rows), columns 1-5.\n#\n# Addressing each cell directly via Table.Cell(row, column) - rather than\n# a document-wide Find/Replace - is important here: several of the new\n# values are identical to *other* cells' old values (e.g. one cell goes\n# from \"95\u00f75=\" to \"26\u00f79=\", while a different cell already holds\n# \"26\u00f79=\" and must become \"93\u00f77=\"). A global find-and-replace-all run\n# sequentially would re-match text that a previous replacement just\n# wrote, corrupting later cells. Targeting the fixed cell coordinates\n# sidesteps that entirely.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$values = @(\n    @(\"34\u00f74=\", \"24\u00f79=\"), @(\"56\u00f73=\", \"18\u00f73=\"), @(\"43\u00f72=\", \"49\u00f73=\"), @(\"95\u00f75=\", \"26\u00f79=\"), @(\"52\u00f75=\", \"99\u00f73=\"),\n    @(\"48\u00f75=\", \"32\u00f76=\"), @(\"12\u00f75=\", \"55\u00f72=\"), @(\"20\u00f73=\", \"62\u00f74=\"), @(\"14\u00f74=\", \"94\u00f79=\"), @(\"65\u00f74=\", \"54\u00f73=\"),\n    @(\"90\u00f77=\", \"96\u00f76=\"), @(\"57\u00f77=\", \"81\u00f77=\"), @(\"80\u00f73=\", \"43\u00f73=\"), @(\"29\u00f75=\", \"57\u00f77=\"), @(\"45\u00f77=\", \"69\u00f72=\"),\n    @(\"19\u00f72=\", \"13\u00f73=\"), @(\"45\u00f79=\", \"38\u00f74=\"), @(\"59\u00f77=\", \"17\u00f79=\"), @(\"26\u00f79=\", \"93\u00f77=\"), @(\"27\u00f77=\", \"49\u00f75=\"),\n    @(\"10\u00f74=\", \"51\u00f76=\"), @(\"74\u00f72=\", \"99\u00f74=\"), @(\"15\u00f76=\", \"23\u00f72=\"), @(\"26\u00f76=\", \"77\u00f77=\"), @(\"57\u00f73=\", \"36\u00f74=\")\n)\n\n$i = 0\nforeach ($row in $dataRows) {\n    for ($col = 1; $col -le 5; $col++) {\n        $pair = $values[$i]\n        $old = $pair[0]\n        $new = $pair[1]\n        $cell = $t.Cell($row, $col)\n        $cellRange = $cell.Range\n        $cellRange.MoveEnd(1, -1) | Out-Null\n        if ($cellRange.Text -ne $old) {\n            throw \"Unexpected text at row $row, col $col. Expected '$old' but found '$($cellRange.Text)'.\"\n        }\n        $cellRange.Text = $new\n        $i++\n    }\n}\n"}
